{"js": "// The document's single \"_GoBack\" bookmark is being moved from the end of\n// the document (right after \"Winter 2019)\") to the title line, right after\n// the first name \"Phillip\" (replacing the old trailing-space run split of\n// \"Phillip \" / \"James \" with \"Phillip\" / \" \").\n\n// 1) Remove the stale _GoBack bookmark whatever/wherever it currently is.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the \"Phillip \" run in the title paragraph (first paragraph of the\n//    document) and trim the trailing space from it.\nconst firstParagraph = context.document.body.paragraphs.getFirst();\nconst nameHits = firstParagraph.search(\"Phillip \", { matchCase: true });\nnameHits.load(\"items\");\nawait context.sync();\n\nconst phillipRange = nameHits.items[0];\nphillipRange.insertText(\"Phillip\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Re-insert the _GoBack bookmark immediately after \"Phillip\".\nconst afterPhillip = phillipRange.getRange(Word.RangeLocation.after);\nafterPhillip.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) The following \"James \" run becomes a single space \" \".\nconst jamesHits = firstParagraph.search(\"James \", { matchCase: true });\njamesHits.load(\"items\");\nawait context.sync();\n\nconst jamesRange = jamesHits.items[0];\njamesRange.insertText(\" \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document's single \"_GoBack\" bookmark is being moved from the end of\n# the document (right after \"Winter 2019)\") to the title line, right after\n# the first name \"Phillip\" (the old trailing-space run split of\n# \"Phillip \" / \"James \" becomes \"Phillip\" / \" \").\n\n$d = $word.ActiveDocument\n\n# 1) Locate the \"Phillip \" run in the title (first) paragraph and collapse\n#    a range to right after it.\n$titleRange = $d.Paragraphs(1).Range.Duplicate\n$titleRange.Find.Execute(\"Phillip \") | Out-Null\n$titleRange.Collapse(0)  # wdCollapseEnd\n\n# 2) Re-add the \"_GoBack\" bookmark at that spot. Bookmark names are unique,\n#    so this automatically removes it from its old location (right after\n#    \"Winter 2019)\") and plants it here instead.\n$d.Bookmarks.Add(\"_GoBack\", $titleRange) | Out-Null\n\n# 3) Trim the trailing space from \"Phillip \".\n$phillipRange = $d.Paragraphs(1).Range.Duplicate\n$phillipRange.Find.Execute(\"Phillip \") | Out-Null\n$phillipRange.Text = \"Phillip\"\n\n# 4) The following \"James \" run becomes a single space \" \".\n$jamesRange = $d.Paragraphs(1).Range.Duplicate\n$jamesRange.Find.Execute(\"James \") | Out-Null\n$jamesRange.Text = \" \"\n"}
